$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'289.61"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'-4.06%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'30.82"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'-4.09%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'4.879"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'-2.17%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.07150"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'-9.53%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'1.852"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'-11.88%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'7.638"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-2.11%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'3.734"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'-1.65%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.8964"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'-3.47%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.1647"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'-6.00%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.07508"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'-5.38%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.08152"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'-5.27%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.02997"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-4.41%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.09996"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-0.16%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.001497"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-1.45%"
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'0.005832"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'0.68%"
$ws.Range("E16").ClearFormats()
$ws.Range("D18").Value = "'3.459"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'-0.10%"
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'-7.40%"
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'0.3276"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'-0.35%"
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'0.1292"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'-1.25%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'4.266"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'-0.27%"
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = "'11.81%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.04473"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'-2.93%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.001213"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'-1.98%"
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'0.004658"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'4.35%"
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'0.13%"
$ws.Range("E27").ClearFormats()
$ws.Range("D39").Value = "'0.01637"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'-4.65%"
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.04340"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'-9.30%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.007442"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-0.32%"
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.1305"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'-4.03%"
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'-16.06%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.01028"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'0.35%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.00005845"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'-2.67%"
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'0.14%"
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'168.35%"
$ws.Range("E47").ClearFormats()
$ws.Range("E48").Value = "'-11.46%"
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.00002104"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'0.14%"
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'0.14%"
$ws.Range("E50").ClearFormats()
